# Galoyan Time recording log: add "Week 2" sheet (copy of "Week 1"),
# clear the template rows, fill in the first log entry, and wrap the
# "Comments" column text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Week 1")

# Duplicate "Week 1" right after itself, then rename the copy to "Week 2".
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Week 2"

# Clear the template's sample data rows (8-18), leaving row indices (A)
# and the Delta Time formulas (F) intact.
$ws2.Range("B8:E18").ClearContents()
$ws2.Range("G8:J18").ClearContents()

# Fill in the first entry of the new week.
$ws2.Range("B7").Value = 43501
$ws2.Range("C7").Value = 0.80208333333333337
$ws2.Range("D7").Value = 0.9375
$ws2.Range("G7").Value = "Prep."
$ws2.Range("H7").Value = "Watching HTML and CSS course in MVA and writing notes"

# Wrap the "Comments" column so the long activity descriptions fit.
$ws2.Range("H7:H18").WrapText = $true
$ws2.Rows(7).RowHeight = 28.8

# Match the saved selection / active sheet state.
$ws2.Range("D8").Select()
